$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$vals = @(78.9,78.92,78.94,78.96,78.98,79.0,79.02,79.04,79.06,79.08,79.10,79.12,79.14,79.16,79.18,79.20,79.22,79.24,79.26,79.28,79.30,79.32,79.34,79.36,79.38,79.40)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Columns.Item($i+1).ColumnWidth = $vals[$i]
}
Write-Output "Done"
